$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 14 was mis-shifted: the title ended up in B14 ("Diagnostic") while the
# real title sat in C14 and the first_author was missing entirely.
# Fix: put the title back in B14, and the first author ("Dewar, LJ") in C14.
$ws.Range("B14").Value = "Investigating the genetic causes of sudden unexpected death in children through targeted next-generation sequencing analysis"
$ws.Range("C14").Value = "Dewar, LJ"

# Remove the now-stale autofilter/sort state left on the sheet.
$ws.AutoFilterMode = $false

# Move the active selection as in the saved workbook.
$ws.Range("M23").Select()
